# Weekly price update: insert two new daily price rows (Ají - Americana (o) / Segunda
# and Primera) for Comercializadora del Agro de Limarí, shifting the existing rows
# (old 68..166) down to (70..168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 68; Excel shifts rows 68:166 down to 70:168
# and copies formatting (e.g. the date number format in column D) from the surrounding rows.
$ws.Rows("68:69").Insert()

# --- New row 68: Ají, Americana (o), Primera ---
$ws.Range("A68").Value = 2
$ws.Range("B68").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 44483
$ws.Range("E68").Value = 4
$ws.Range("F68").Value = 100112021
$ws.Range("G68").Value = "Ají"
$ws.Range("H68").Value = "Americana (o)"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 100
$ws.Range("K68").Value = 50000
$ws.Range("L68").Value = 55000
$ws.Range("M68").Value = 52500
$ws.Range("N68").Value = "$/caja 25 kilos"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 2100
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"

# --- New row 69: Ají, Americana (o), Segunda ---
$ws.Range("A69").Value = 2
$ws.Range("B69").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 44483
$ws.Range("E69").Value = 4
$ws.Range("F69").Value = 100112021
$ws.Range("G69").Value = "Ají"
$ws.Range("H69").Value = "Americana (o)"
$ws.Range("I69").Value = "Segunda"
$ws.Range("J69").Value = 40
$ws.Range("K69").Value = 40000
$ws.Range("L69").Value = 45000
$ws.Range("M69").Value = 42500
$ws.Range("N69").Value = "$/caja 25 kilos"
$ws.Range("O69").Value = "Provincia de Limarí"
$ws.Range("P69").Value = 1700
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"
